$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff Datetime (D5) and Correspond Handback DateTime (G5)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-25 11:21:21"
$wsZhCn.Range("G5").Value = "2016-01-25 11:22:02"

# de-de sheet: update Correspond Handoff Datetime (D5) and Correspond Handback DateTime (G5)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-25 11:21:30"
$wsDeDe.Range("G5").Value = "2016-01-25 11:22:19"
